$wb = $excel.ActiveWorkbook

# ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1015.875
$ws.Range("I43").Value = 650
$ws.Range("J43").Value = 1137.8334
$ws.Range("K43").Value = 650
$ws.Range("L43").Value = 1137.8334
$ws.Range("M43").Value = -581
$ws.Range("N43").Value = -1275.8334

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2225.375
$ws.Range("I62").Value = 2225.375
$ws.Range("K62").Value = 2225.375
$ws.Range("M62").Value = -1601.375

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2225.375
$ws.Range("I65").Value = 2225.375
$ws.Range("K65").Value = 11126.875
$ws.Range("M65").Value = -8006.875

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1066.459
$ws.Range("J129").Value = 1106
$ws.Range("L129").Value = 3318
$ws.Range("N129").Value = -13318

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1300.1621
$ws.Range("I132").Value = 1065.6875
$ws.Range("K132").Value = 3197.0625
$ws.Range("M132").Value = -667.0625

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 3393.5854
$ws.Range("I135").Value = 2028.2307
$ws.Range("J135").Value = 30018
$ws.Range("K135").Value = 18254.0763
$ws.Range("L135").Value = 270162
$ws.Range("M135").Value = -15719.0763
$ws.Range("N135").Value = -275232

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1955.1111
$ws.Range("I137").Value = 1854.8889
$ws.Range("J137").Value = 2155.5557
$ws.Range("K137").Value = 5564.6667
$ws.Range("L137").Value = 6466.6671
$ws.Range("M137").Value = -3014.6667
$ws.Range("N137").Value = -11566.6671

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2483.75
$ws.Range("I141").Value = 2042.1875
$ws.Range("J141").Value = 4250
$ws.Range("K141").Value = 6126.5625
$ws.Range("L141").Value = 12750
$ws.Range("M141").Value = -946.5625
$ws.Range("N141").Value = -23110

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7552.578
$ws.Range("I32").Value = 5826.9043
$ws.Range("K32").Value = 5826.9043
$ws.Range("M32").Value = -5539.9043

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 18262
$ws.Range("I45").Value = 21314.4
$ws.Range("K45").Value = 21314.4
$ws.Range("M45").Value = -20937.4

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2110.9412
$ws.Range("I74").Value = 2250
$ws.Range("J74").Value = 1954.5
$ws.Range("K74").Value = 2250
$ws.Range("L74").Value = 1954.5
$ws.Range("M74").Value = -1376
$ws.Range("N74").Value = -3702.5

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2110.9412
$ws.Range("I77").Value = 2250
$ws.Range("J77").Value = 1954.5
$ws.Range("K77").Value = 11250
$ws.Range("L77").Value = 9772.5
$ws.Range("M77").Value = -6882
$ws.Range("N77").Value = -18508.5

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1123.5714
$ws.Range("I97").Value = 734.3333
$ws.Range("J97").Value = 1824.2
$ws.Range("K97").Value = 734.3333
$ws.Range("L97").Value = 1824.2
$ws.Range("M97").Value = -238.3333
$ws.Range("N97").Value = -2816.2

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1268.6666
$ws.Range("I110").Value = 1299.8
$ws.Range("K110").Value = 1299.8
$ws.Range("M110").Value = 745.2

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2422.9375
$ws.Range("I94").Value = 1851
$ws.Range("J94").Value = 2994.875
$ws.Range("K94").Value = 1851
$ws.Range("L94").Value = 2994.875
$ws.Range("M94").Value = -1400
$ws.Range("N94").Value = -3896.875

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1720.409
$ws.Range("I16").Value = 1434
$ws.Range("J16").Value = 2334.1428
$ws.Range("K16").Value = 1434
$ws.Range("L16").Value = 2334.1428
$ws.Range("M16").Value = -1147
$ws.Range("N16").Value = -2908.1428

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 899.6
$ws.Range("I22").Value = 833
$ws.Range("J22").Value = 999.5
$ws.Range("K22").Value = 833
$ws.Range("L22").Value = 999.5
$ws.Range("M22").Value = -483
$ws.Range("N22").Value = -1699.5

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3448.9893
$ws.Range("I31").Value = 1390.902
$ws.Range("J31").Value = 5948.095
$ws.Range("K31").Value = 1390.902
$ws.Range("L31").Value = 5948.095
$ws.Range("M31").Value = -1095.902
$ws.Range("N31").Value = -6538.095

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3448.9893
$ws.Range("I34").Value = 1390.902
$ws.Range("J34").Value = 5948.095
$ws.Range("K34").Value = 1390.902
$ws.Range("L34").Value = 5948.095
$ws.Range("M34").Value = -1188.902
$ws.Range("N34").Value = -6352.095

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2458.5625
$ws.Range("I58").Value = 1953.8334
$ws.Range("J58").Value = 2761.4
$ws.Range("K58").Value = 1953.8334
$ws.Range("L58").Value = 2761.4
$ws.Range("M58").Value = -1750.8334
$ws.Range("N58").Value = -3167.4

# CRP row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 4050.3667
$ws.Range("I94").Value = 3112.7144
$ws.Range("J94").Value = 4870.8125
$ws.Range("K94").Value = 3112.7144
$ws.Range("L94").Value = 4870.8125
$ws.Range("M94").Value = -2661.7144
$ws.Range("N94").Value = -5772.8125

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3307.96
$ws.Range("I99").Value = 3685.6667
$ws.Range("J99").Value = 1325
$ws.Range("K99").Value = 3685.6667
$ws.Range("L99").Value = 1325
$ws.Range("M99").Value = -2187.6667
$ws.Range("N99").Value = -4321

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1808.7646
$ws.Range("I105").Value = 1795.75
$ws.Range("J105").Value = 1840
$ws.Range("K105").Value = 1795.75
$ws.Range("L105").Value = 1840
$ws.Range("M105").Value = -48.75
$ws.Range("N105").Value = -5334

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1720.409
$ws.Range("I113").Value = 1434
$ws.Range("J113").Value = 2334.1428
$ws.Range("K113").Value = 1434
$ws.Range("L113").Value = 2334.1428
$ws.Range("M113").Value = 736
$ws.Range("N113").Value = -6674.1428

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2317895.8
$ws.Range("I122").Value = 3968670.5
$ws.Range("J122").Value = 6811.2
$ws.Range("K122").Value = 11906011.5
$ws.Range("L122").Value = 20433.6
$ws.Range("M122").Value = -11903561.5
$ws.Range("N122").Value = -25333.6

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3307.96
$ws.Range("I126").Value = 3685.6667
$ws.Range("J126").Value = 1325
$ws.Range("K126").Value = 11057.0001
$ws.Range("L126").Value = 3975
$ws.Range("M126").Value = -8587.000100000001
$ws.Range("N126").Value = -8915

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2051.0688
$ws.Range("I132").Value = 1662.7727
$ws.Range("J132").Value = 3271.4285
$ws.Range("K132").Value = 4988.3181
$ws.Range("L132").Value = 9814.2855
$ws.Range("M132").Value = -2458.3181
$ws.Range("N132").Value = -14874.2855

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2217.9805
$ws.Range("I134").Value = 2272.976
$ws.Range("J134").Value = 1961.3334
$ws.Range("K134").Value = 6818.928
$ws.Range("L134").Value = 5884.0002
$ws.Range("M134").Value = -4283.928
$ws.Range("N134").Value = -10954.0002

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2458.5625
$ws.Range("I136").Value = 1953.8334
$ws.Range("J136").Value = 2761.4
$ws.Range("K136").Value = 5861.5002
$ws.Range("L136").Value = 8284.200000000001
$ws.Range("M136").Value = -3311.5002
$ws.Range("N136").Value = -13384.2

# CUL row 23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 8333463
$ws.Range("I23").Value = 20000090
$ws.Range("J23").Value = 157.42857
$ws.Range("K23").Value = 60000270
$ws.Range("L23").Value = 472.28571
$ws.Range("M23").Value = -60000035
$ws.Range("N23").Value = -942.28571

# CUL row 69
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 1348.5714
$ws.Range("I69").Value = 637.3333
$ws.Range("J69").Value = 1882
$ws.Range("K69").Value = 1911.9999
$ws.Range("L69").Value = 5646
$ws.Range("M69").Value = -1100.9999
$ws.Range("N69").Value = -7268

# CUL row 72
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 1348.5714
$ws.Range("I72").Value = 637.3333
$ws.Range("J72").Value = 1882
$ws.Range("K72").Value = 5735.9997
$ws.Range("L72").Value = 16938
$ws.Range("M72").Value = -1679.9997
$ws.Range("N72").Value = -25050

# CUL row 117
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 17550878
$ws.Range("I117").Value = 25407.25
$ws.Range("J117").Value = 22224336
$ws.Range("K117").Value = 76221.75
$ws.Range("L117").Value = 66673008
$ws.Range("M117").Value = -72779.75
$ws.Range("N117").Value = -66679892

# CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 911.58826
$ws.Range("I129").Value = 721.5
$ws.Range("J129").Value = 1015.2727
$ws.Range("K129").Value = 2164.5
$ws.Range("L129").Value = 3045.8181
$ws.Range("M129").Value = 2835.5
$ws.Range("N129").Value = -13045.8181

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5446.467
$ws.Range("I70").Value = 5412.8696
$ws.Range("J70").Value = 5556.857
$ws.Range("K70").Value = 5412.8696
$ws.Range("L70").Value = 5556.857
$ws.Range("M70").Value = -5142.8696
$ws.Range("N70").Value = -6096.857

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5446.467
$ws.Range("I73").Value = 5412.8696
$ws.Range("J73").Value = 5556.857
$ws.Range("K73").Value = 5412.8696
$ws.Range("L73").Value = 5556.857
$ws.Range("M73").Value = -4476.8696
$ws.Range("N73").Value = -7428.857

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 50001380
$ws.Range("I113").Value = 62500924
$ws.Range("J113").Value = 3200
$ws.Range("K113").Value = 62500924
$ws.Range("L113").Value = 3200
$ws.Range("M113").Value = -62498754
$ws.Range("N113").Value = -7540

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5540.8335
$ws.Range("I126").Value = 10034.083
$ws.Range("J126").Value = 2545.3333
$ws.Range("K126").Value = 30102.249
$ws.Range("L126").Value = 7635.999899999999
$ws.Range("M126").Value = -27632.249
$ws.Range("N126").Value = -12575.9999

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 174.36111
$ws.Range("I55").Value = 164.55556
$ws.Range("J55").Value = 184.16667
$ws.Range("K55").Value = 164.55556
$ws.Range("L55").Value = 184.16667
$ws.Range("M55").Value = 8.444439999999986
$ws.Range("N55").Value = -530.1666700000001

# WVR row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 8257.143
$ws.Range("J54").Value = 8257.143
$ws.Range("L54").Value = 8257.143
$ws.Range("N54").Value = -9297.143

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1421.5714
$ws.Range("I81").Value = 1421.5714
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2843.1428
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -1782.1428
$ws.Range("N81").ClearContents()

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1421.5714
$ws.Range("I84").Value = 1421.5714
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 14215.714
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -8911.714
$ws.Range("N84").ClearContents()
